# Apply the edits described by the diff to the active workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the worksheet's display name (tab stays same sheetId/rId).
$ws.Name = "CubeA"

# 2) Tiny floating point corrections in existing rows (recalculation artifacts).
$ws.Range("M13").Value = 0.9943642007518049
$ws.Range("N13").Value = 1.008226233984116
$ws.Range("J15").Value = 0.9047947004779515
$ws.Range("L15").Value = 0.6924514826557866
$ws.Range("P15").Value = 1.001718261182721

# 3) Append a new data row (row 16) with the Gaussian Quadrature results.
#    Copy the formatting from the row above (A15 carries the bold/bordered style)
#    so the new index cell matches the rest of column A.
$ws.Range("A15").Copy($ws.Range("A16"))
$ws.Range("A16").Value = 14

$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"

$ws.Range("C16").Value = 0.2971667452826968
$ws.Range("D16").Value = 4.560219924018299
$ws.Range("E16").Value = 2.001981095267751
$ws.Range("F16").Value = 0.6088131020872382
$ws.Range("G16").Value = 0.2971667452826968
$ws.Range("H16").Value = 4.560219924018299
$ws.Range("I16").Value = 1.021121545906042
$ws.Range("J16").Value = 1.001612662612057
$ws.Range("K16").Value = 0.5187705011526013
$ws.Range("L16").Value = 1.183750921905119
$ws.Range("M16").Value = 0.2971667452826968
$ws.Range("N16").Value = 3.281100509643025
$ws.Range("O16").Value = 1.867045216663996
$ws.Range("P16").Value = 1.399179562278976
